# Auto-generated edit script applying numeric updates described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3541.8667
$ws.Range("I40").Value = 1037
$ws.Range("K40").Value = 1037
$ws.Range("M40").Value = -862

$ws.Range("H53").Value = 875.4737
$ws.Range("I53").Value = 104.6
$ws.Range("J53").Value = 1150.7858
$ws.Range("K53").Value = 104.6
$ws.Range("L53").Value = 1150.7858
$ws.Range("M53").Value = 532.4
$ws.Range("N53").Value = -2424.7858

$ws.Range("H55").Value = 419.57144
$ws.Range("I55").Value = 59
$ws.Range("K55").Value = 59
$ws.Range("M55").Value = 155

$ws.Range("H86").Value = 6079.5
$ws.Range("I86").Value = 1424
$ws.Range("J86").Value = 7010.6
$ws.Range("K86").Value = 1424
$ws.Range("L86").Value = 7010.6
$ws.Range("M86").Value = -301
$ws.Range("N86").Value = -9256.6

$ws.Range("H89").Value = 6079.5
$ws.Range("I89").Value = 1424
$ws.Range("J89").Value = 7010.6
$ws.Range("K89").Value = 7120
$ws.Range("L89").Value = 35053
$ws.Range("M89").Value = -1504
$ws.Range("N89").Value = -46285

$ws.Range("H115").Value = 4972.625
$ws.Range("I115").Value = 3125.8572
$ws.Range("J115").Value = 17900
$ws.Range("K115").Value = 9377.571599999999
$ws.Range("L115").Value = 53700
$ws.Range("M115").Value = -7810.571599999999
$ws.Range("N115").Value = -56834

$ws.Range("H132").Value = 8591.487999999999
$ws.Range("I132").Value = 3159
$ws.Range("K132").Value = 9477
$ws.Range("M132").Value = -6947

$ws.Range("H133").Value = 101798
$ws.Range("J133").Value = 101798
$ws.Range("L133").Value = 101798
$ws.Range("N133").Value = -111918

$ws.Range("H134").Value = 68612.25
$ws.Range("J134").Value = 68612.25
$ws.Range("L134").Value = 68612.25
$ws.Range("N134").Value = -78752.25

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws.Range("H138").Value = 356067.16
$ws.Range("J138").Value = 1112946.8
$ws.Range("L138").Value = 3338840.4
$ws.Range("N138").Value = -3349120.4

$ws.Range("H139").Value = 85000
$ws.Range("J139").Value = 85000
$ws.Range("L139").Value = 85000
$ws.Range("N139").Value = -95280

$ws.Range("H141").Value = 3104.3684
$ws.Range("J141").Value = 5832
$ws.Range("L141").Value = 17496
$ws.Range("N141").Value = -27856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3464.5
$ws.Range("I45").Value = 3126.8333
$ws.Range("K45").Value = 3126.8333
$ws.Range("M45").Value = -2749.8333

$ws.Range("H61").Value = 8047.8945
$ws.Range("I61").Value = 4551.25
$ws.Range("J61").Value = 10590.909
$ws.Range("K61").Value = 4551.25
$ws.Range("L61").Value = 10590.909
$ws.Range("M61").Value = -4339.25
$ws.Range("N61").Value = -11014.909

$ws.Range("H97").Value = 1406.6666
$ws.Range("I97").Value = 1406.6666
$ws.Range("K97").Value = 1406.6666
$ws.Range("M97").Value = -910.6666

$ws.Range("H110").Value = 4261.5
$ws.Range("I110").Value = 2800
$ws.Range("J110").Value = 4992.25
$ws.Range("K110").Value = 2800
$ws.Range("L110").Value = 4992.25
$ws.Range("M110").Value = -755
$ws.Range("N110").Value = -9082.25

$ws.Range("H132").Value = 2151.6775
$ws.Range("I132").Value = 1627.1111
$ws.Range("J132").Value = 5692.5
$ws.Range("K132").Value = 4881.3333
$ws.Range("L132").Value = 17077.5
$ws.Range("M132").Value = -2351.3333
$ws.Range("N132").Value = -22137.5

$ws.Range("H136").Value = 8047.8945
$ws.Range("I136").Value = 4551.25
$ws.Range("J136").Value = 10590.909
$ws.Range("K136").Value = 13653.75
$ws.Range("L136").Value = 31772.727
$ws.Range("M136").Value = -11103.75
$ws.Range("N136").Value = -36872.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 637.61536
$ws.Range("I22").Value = 657.4167
$ws.Range("K22").Value = 657.4167
$ws.Range("M22").Value = -484.4167

$ws.Range("H94").Value = 1098.1428
$ws.Range("I94").Value = 234.25
$ws.Range("J94").Value = 2250
$ws.Range("K94").Value = 234.25
$ws.Range("L94").Value = 2250
$ws.Range("M94").Value = 216.75
$ws.Range("N94").Value = -3152

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2658.577
$ws.Range("J31").Value = 7699.6
$ws.Range("L31").Value = 7699.6
$ws.Range("N31").Value = -8289.6

$ws.Range("H34").Value = 2658.577
$ws.Range("J34").Value = 7699.6
$ws.Range("L34").Value = 7699.6
$ws.Range("N34").Value = -8103.6

$ws.Range("H58").Value = 2570.125
$ws.Range("I58").Value = 3266.75
$ws.Range("K58").Value = 3266.75
$ws.Range("M58").Value = -3063.75

$ws.Range("H107").Value = 1316.5416
$ws.Range("I107").Value = 485.27274
$ws.Range("K107").Value = 485.27274
$ws.Range("M107").Value = 1434.72726

$ws.Range("H132").Value = 1540927.6
$ws.Range("I132").Value = 1907205.8
$ws.Range("K132").Value = 5721617.4
$ws.Range("M132").Value = -5719087.4

$ws.Range("H134").Value = 2272.5095
$ws.Range("I134").Value = 1219.6046
$ws.Range("K134").Value = 3658.8138
$ws.Range("M134").Value = -1123.8138

$ws.Range("H136").Value = 2570.125
$ws.Range("I136").Value = 3266.75
$ws.Range("K136").Value = 9800.25
$ws.Range("M136").Value = -7250.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2555.5833
$ws.Range("J75").Value = 3436.4285
$ws.Range("L75").Value = 10309.2855
$ws.Range("N75").Value = -12305.2855

$ws.Range("H78").Value = 2555.5833
$ws.Range("J78").Value = 3436.4285
$ws.Range("L78").Value = 30927.8565
$ws.Range("N78").Value = -40911.8565

$ws.Range("H92").Value = 740.4286
$ws.Range("I92").Value = 291.33334
$ws.Range("J92").Value = 862.9091
$ws.Range("K92").Value = 874.0000200000001
$ws.Range("L92").Value = 2588.7273
$ws.Range("M92").Value = 373.9999799999999
$ws.Range("N92").Value = -5084.7273

$ws.Range("H103").Value = 306.7
$ws.Range("I103").Value = 118.75
$ws.Range("J103").Value = 432
$ws.Range("K103").Value = 356.25
$ws.Range("L103").Value = 1296
$ws.Range("M103").Value = 522.75
$ws.Range("N103").Value = -3054

$ws.Range("H141").Value = 8006.846
$ws.Range("I141").Value = 6006.364
$ws.Range("J141").Value = 19009.5
$ws.Range("K141").Value = 18019.092
$ws.Range("L141").Value = 57028.5
$ws.Range("M141").Value = -12839.092
$ws.Range("N141").Value = -67388.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5049.3687
$ws.Range("I70").Value = 4595.273
$ws.Range("K70").Value = 4595.273
$ws.Range("M70").Value = -4325.273

$ws.Range("H73").Value = 5049.3687
$ws.Range("I73").Value = 4595.273
$ws.Range("K73").Value = 4595.273
$ws.Range("M73").Value = -3659.273

$ws.Range("H97").Value = 3359.8333
$ws.Range("I97").Value = 2832
$ws.Range("J97").Value = 5999
$ws.Range("K97").Value = 2832
$ws.Range("L97").Value = 5999
$ws.Range("M97").Value = -2336
$ws.Range("N97").Value = -6991

$ws.Range("H102").Value = 36778.62
$ws.Range("I102").Value = 2201.7144
$ws.Range("K102").Value = 2201.7144
$ws.Range("M102").Value = -579.7143999999998

$ws.Range("H122").Value = 2611.389
$ws.Range("I122").Value = 2371.7856
$ws.Range("J122").Value = 3450
$ws.Range("K122").Value = 7115.3568
$ws.Range("L122").Value = 10350
$ws.Range("M122").Value = -4665.3568
$ws.Range("N122").Value = -15250

$ws.Range("H126").Value = 1996.6666
$ws.Range("I126").Value = 1495
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 4485
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -2015
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 13892035
$ws.Range("I132").Value = 14495949
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 43487847
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -43485317
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8701505
$ws.Range("I7").Value = 15389009
$ws.Range("K7").Value = 15389009
$ws.Range("M7").Value = -15388897

$ws.Range("H55").Value = 168.85715
$ws.Range("I55").Value = 95.70587999999999
$ws.Range("K55").Value = 95.70587999999999
$ws.Range("M55").Value = 77.29412000000001

$ws.Range("H122").Value = 4289.8
$ws.Range("I122").Value = 2500.2
$ws.Range("J122").Value = 6079.4
$ws.Range("K122").Value = 7500.599999999999
$ws.Range("L122").Value = 18238.2
$ws.Range("M122").Value = -5050.599999999999
$ws.Range("N122").Value = -23138.2

$ws.Range("H126").Value = 8701505
$ws.Range("I126").Value = 15389009
$ws.Range("K126").Value = 46167027
$ws.Range("M126").Value = -46164557

$ws.Range("H132").Value = 3131.432
$ws.Range("I132").Value = 3150.625
$ws.Range("J132").Value = 3080.25
$ws.Range("K132").Value = 9451.875
$ws.Range("L132").Value = 9240.75
$ws.Range("M132").Value = -6921.875
$ws.Range("N132").Value = -14300.75

$ws.Range("H136").Value = 4188.148
$ws.Range("I136").Value = 3952.2104
$ws.Range("J136").Value = 4748.5
$ws.Range("K136").Value = 11856.6312
$ws.Range("L136").Value = 14245.5
$ws.Range("M136").Value = -9306.6312
$ws.Range("N136").Value = -19345.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 21302.25
$ws.Range("I51").Value = 21302.25
$ws.Range("K51").Value = 21302.25
$ws.Range("M51").Value = -20792.25

$ws.Range("H122").Value = 5355.8057
$ws.Range("I122").Value = 5424.9697
$ws.Range("K122").Value = 16274.9091
$ws.Range("M122").Value = -13824.9091

$ws.Range("H132").Value = 1652.4762
$ws.Range("I132").Value = 1554.7736
$ws.Range("K132").Value = 4664.3208
$ws.Range("M132").Value = -2134.3208

$ws.Range("H136").Value = 20153.637
$ws.Range("J136").Value = 997.1429000000001
$ws.Range("L136").Value = 2991.4287
$ws.Range("N136").Value = -8091.4287

